$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.618.10"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "3.298.00"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "3.290.44"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.571"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("E10").Value = "  -6.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.572"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").Value = "3.828.89"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "568.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -10.48%  "
$ws.Range("D17").Value = "65.562.01"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "3.299.16"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("E21").Value = "  -5.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.882"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.97%  "
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "554.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("D35").Value = "3.750.63"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("E41").Value = "  -8.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "0.0₃0667"
$ws.Range("E43").Value = "  -8.34%  "
$ws.Range("E44").Value = "  -6.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.327"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.33%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.00%  "
